$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: "Entrega 3" indicators, mirroring column B's formatting ---
# Seed C1:C4 styles by copying the matching B-column cell's format, then
# overwrite with the new values/text.
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("B2").Copy($ws.Range("C2"))
$ws.Range("B3").Copy($ws.Range("C3"))
$ws.Range("B4").Copy($ws.Range("C4"))

# Header row: B1 becomes "Entrega 2", C1 becomes "Entrega 3"
$ws.Range("B1").Value2 = "Indicadores de Qualidade - Entrega 2"
$ws.Range("C1").Value2 = "Indicadores de Qualidade - Entrega 3"

# New data column values
$ws.Range("C2").Value2 = 4
$ws.Range("C3").Value2 = 4
$ws.Range("C4").Value2 = 5

# Header row gets taller to fit the wrapped two-line titles
$ws.Range("A1:C1").RowHeight = 25.5

# Trailing blank row/cell (A13) present but empty, with the default style
$ws.Range("A13").Value2 = "x"
$ws.Range("A13").ClearContents()
$ws.Range("A13").Style = "Normal"

# Match the committed selection state
$ws.Range("A1:C12").Select() | Out-Null
